$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove row 86 ("La misère est si belle (MAYKO Remix)") entirely; subsequent
# rows (87-90) shift up to become rows 86-89, and the used range shrinks
# from A1:C90 to A1:C89.
$ws.Rows.Item(86).Delete()
